# Append 4 new applicant rows (19-22) to Sheet1, mirroring the existing
# data rows which are all stored as text (inlineStr) cells - including
# the numeric-looking Phone Number / Experience columns.
#
# A plain `$cell.Value = "8798569826"` would be auto-coerced to a Number
# by Excel's smart-typing, so instead we build each value with a tiny
# text-returning formula in a scratch cell, copy it, and paste-special
# just the (already-text-typed) value into the destination cell. That
# keeps the destination a plain text cell with no extra number format /
# style, matching how the rest of the sheet is stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Cells.Item(200, 1)

function Set-TextCell($row, $col, $text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
}

$rows = @(
    @("Praveen", "Masters", "praveen@gmail.com", "8798569826", "5", "pravven"),
    @("Ram", "M.Sc", "ram@gmail.com", "988765487", "2", "78945"),
    @("Mani", "M.Sc", "praveen@gmail.com", "78932145610", "0", "mani@44"),
    @("Subramanian S", "M. Sc", "subramanian160104@gmail.com", "8072744511", "1", "123mani")
)

$startRow = 19
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowValues = $rows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $col = $c + 1
        Set-TextCell $r $col $rowValues[$c]
    }
}

$scratch.Clear()
